# Refresh quarterly macro-data download with the latest source pull.
# Updates a handful of already-populated cells (Population interpolation in
# column H, Unemployment Rate in column O, a couple of one-off corrections)
# and fills in newly-available figures for the most recent quarters
# (rows 116-120), matching the newly released
# Inflation_contributions_graph_data.XLSX / Waterfall_graph_data.XLSX /
# trimmed_graph_data.XLSX / MEGA_DATA_DOWNLOAD.xlsx /
# SUMMARY_EXPORT_DATA_DOWNLOAD.xlsx / EXPORT_DATA_DOWNLOAD_ALL.xlsx sources.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unemployment Rate (column O) historical restatements ---
$ws.Range("O14").Value = 8.83333333333333
$ws.Range("O16").Value = 7.26666666666667
$ws.Range("O20").Value = 5.43333333333333

# --- Population (column H) interpolated series refresh, rows 66-119 ---
$ws.Range("H66").Value = 4575
$ws.Range("H67").Value = 4580
$ws.Range("H68").Value = 4584
$ws.Range("H69").Value = 4589
$ws.Range("H70").Value = 4593.7
$ws.Range("H71").Value = 4599
$ws.Range("H72").Value = 4604
$ws.Range("O72").Value = 15.4333333333333
$ws.Range("H73").Value = 4609
$ws.Range("H74").Value = 4614.7
$ws.Range("O74").Value = 14.5333333333333
$ws.Range("H75").Value = 4622
$ws.Range("H76").Value = 4630
$ws.Range("H77").Value = 4638
$ws.Range("H78").Value = 4645.4
$ws.Range("H79").Value = 4656
$ws.Range("H80").Value = 4667
$ws.Range("H81").Value = 4677
$ws.Range("H82").Value = 4687.8
$ws.Range("H83").Value = 4701
$ws.Range("O83").Value = 10.1
$ws.Range("H84").Value = 4714
$ws.Range("H85").Value = 4727
$ws.Range("H86").Value = 4739.6
$ws.Range("H87").Value = 4757
$ws.Range("H88").Value = 4775
$ws.Range("H89").Value = 4793
$ws.Range("H90").Value = 4810.9
$ws.Range("H91").Value = 4829
$ws.Range("H92").Value = 4848
$ws.Range("O92").Value = 6.66666666666667
$ws.Range("H93").Value = 4866
$ws.Range("H94").Value = 4884.9
$ws.Range("H95").Value = 4903
$ws.Range("H96").Value = 4922
$ws.Range("H97").Value = 4940
$ws.Range("O97").Value = 5.76666666666667
$ws.Range("H98").Value = 4958.5
$ws.Range("H99").Value = 4976
$ws.Range("H100").Value = 4994
$ws.Range("H101").Value = 5012
$ws.Range("H102").Value = 5029.9
$ws.Range("H103").Value = 5041
$ws.Range("H104").Value = 5052
$ws.Range("H105").Value = 5064
$ws.Range("H106").Value = 5074.7
$ws.Range("O106").Value = 7.36666666666667
$ws.Range("H107").Value = 5102.025
$ws.Range("O107").Value = 6.96666666666667
$ws.Range("H108").Value = 5129.35
$ws.Range("H109").Value = 5156.675
$ws.Range("H110").Value = 5184
$ws.Range("O110").Value = 4.86666666666667
$ws.Range("H111").Value = 5208.4
$ws.Range("H112").Value = 5232.8
$ws.Range("O112").Value = 4.26666666666667
$ws.Range("H113").Value = 5257.2
$ws.Range("H114").Value = 5281.6
$ws.Range("H115").Value = 5306.275
$ws.Range("O115").Value = 4.16666666666667

# --- Rows 116-119: Population refresh + newly-available Household Deposits (J) ---
$ws.Range("H116").Value = 5330.95
$ws.Range("J116").Value = 456388
$ws.Range("O116").Value = 4.43333333333333

$ws.Range("H117").Value = 5355.625
$ws.Range("J117").Value = 456618

$ws.Range("H118").Value = 5
$ws.Range("J118").Value = 460207
$ws.Range("O118").Value = 4.3
$ws.Range("R118").Value = 731.954

$ws.Range("H119").Value = 380
$ws.Range("J119").Value = 466376
$ws.Range("N119").Value = 2492100
$ws.Range("O119").Value = 4.4

# --- Row 120: newly-available quarter populated across the board ---
$ws.Range("C120").Value = 6988518
$ws.Range("D120").Value = 108
$ws.Range("E120").Value = 21469267453.5
$ws.Range("F120").Value = 157
$ws.Range("J120").Value = 471580
$ws.Range("M120").Value = 2794.8
$ws.Range("N120").Value = 2501633.33333333
$ws.Range("O120").Value = 4.23333333333333
